$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the naive forecaster bug removed the stray y_0_forecast value for
# 2008->2009 (column C) and nudged the recomputed y_1_forecast value.
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 10.06916370210016

# Row 3
$ws.Range("C3").Value = -20.40984652067477
$ws.Range("E3").Value = -33.45158075171639

# Row 4
$ws.Range("C4").Value = 5.331710924091815

# Row 5
$ws.Range("C5").Value = 11.04982736891555
$ws.Range("E5").Value = 9.5899211611429

# Row 6
$ws.Range("E6").Value = 7.086193663491014

# Row 7
$ws.Range("C7").Value = -2.313034291448768

# Row 9
$ws.Range("C9").Value = 3.860244074450203
$ws.Range("E9").Value = 3.254220449867051

# Row 11
$ws.Range("C11").Value = 4.421855465610269
$ws.Range("E11").Value = 5.26036486209962

# Row 12
$ws.Range("C12").Value = 3.320585727896574
$ws.Range("E12").Value = -4.09821347263859

# Row 13
$ws.Range("C13").Value = 1.782333336406383

# Row 14
$ws.Range("E14").Value = -11.8352240479

# Row 15
$ws.Range("C15").Value = 6.928818429977723

# Row 17
$ws.Range("C17").Value = 0.8562564928550342
$ws.Range("E17").Value = 1.697198638953612

# Row 18
$ws.Range("C18").Value = -0.1644433828108638
